# Apply the "Inclusão URL Git ppt" edit to slide 1:
#   1) Split the existing "F1rst Tecnologia" run into "F1rst " + "Tecnologia"
#   2) Add a new text box with the project's Git URL

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# -----------------------------------------------------------------
# 1) Re-split "F1rst Tecnologia" into two runs with identical rPr
# -----------------------------------------------------------------
$subtitle = $s.Shapes.Item(2)
$subTr = $subtitle.TextFrame.TextRange
$subTr.Delete()
$subTr.Text = "F1rst "
$subTr.InsertAfter("Tecnologia")

# -----------------------------------------------------------------
# 2) Add the new "Git: <url>" text box
# -----------------------------------------------------------------
$left   = 214282  / 12700
$top    = 6143644 / 12700
$width  = 5715040 / 12700
$height = 830997  / 12700

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "CaixaDeTexto 3"
$tb.Fill.Visible = $false

$tf = $tb.TextFrame
$tf.WordWrap = -1
$tf.AutoSize = 1

$tr = $tf.TextRange
$tr.Delete()

# Paragraph 1: "Git" + ": " (bold, 16pt, pt-BR)
$tr.Text = "Git: "
$tr.Font.Size = 16
$tr.Font.Bold = $true
$tr.LanguageID = "pt-BR"
$tr.Characters(1, 3).Font.Size = 16
$tr.Characters(4, 2).Font.Size = 16
$tr.InsertAfter([char]13)

# Paragraph 2: "https" + "://" + "github.com/renatomachadosoares/rms-data-master"
$full = $tf.TextRange
$full.InsertAfter("https://github.com/renatomachadosoares/rms-data-master")

$full = $tf.TextRange
$urlRange = $full.Characters(7, 54)
$urlRange.Font.Size = 16
$urlRange.Font.Bold = $false
$urlRange.LanguageID = "pt-BR"
$full.Characters(7, 5).Font.Size = 16
$full.Characters(12, 3).Font.Size = 16
$full.Characters(15, 46).Font.Size = 16

# Paragraph 3: empty trailing paragraph
$full = $tf.TextRange
$full.InsertAfter([char]13)

$subText = $subtitle.TextFrame.TextRange.Text
$boxText = $tf.TextRange.Text
$l = $tb.Left*12700
$t = $tb.Top*12700
$w = $tb.Width*12700
$h = $tb.Height*12700
Write-Host "Subtitle text: [$subText]"
Write-Host "Textbox text: [$boxText]"
Write-Host "Textbox pos/size (EMU): $l $t $w $h"
